$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 183 - odds refresh
$ws.Range("R183").Value2 = 1.87
$ws.Range("S183").Value2 = 2.03
$ws.Range("U183").Value2 = 1.925
$ws.Range("V183").Value2 = 1.925

# Row 185 - now holds match that used to be on row 187 (Sarpsborg vs Odd BK), with refreshed odds
$ws.Range("B185").Value2 = 7617328
$ws.Range("F185").Value2 = "Sarpsborg"
$ws.Range("G185").Value2 = "Odd BK"
$ws.Range("K185").Value2 = 1.571
$ws.Range("L185").Value2 = 4.333
$ws.Range("M185").Value2 = 5
$ws.Range("N185").Value2 = 1.533
$ws.Range("O185").Value2 = 4.5
$ws.Range("P185").Value2 = 5.25
$ws.Range("Q185").Value2 = -1
$ws.Range("R185").Value2 = 1.88
$ws.Range("S185").Value2 = 2.02
$ws.Range("T185").Value2 = 3.25
$ws.Range("U185").Value2 = 1.9
$ws.Range("V185").Value2 = 1.95

# Row 186 - now holds match that used to be on row 185 (Haugesund vs Lillestrom), with refreshed odds
$ws.Range("B186").Value2 = 7617326
$ws.Range("F186").Value2 = "Haugesund"
$ws.Range("G186").Value2 = "Lillestrom"
$ws.Range("K186").Value2 = 2.4
$ws.Range("L186").Value2 = 3.5
$ws.Range("M186").Value2 = 2.75
$ws.Range("N186").Value2 = 2.625
$ws.Range("O186").Value2 = 3.6
$ws.Range("P186").Value2 = 2.5
$ws.Range("Q186").Value2 = 0
$ws.Range("R186").Value2 = 2.02
$ws.Range("S186").Value2 = 1.88
$ws.Range("T186").Value2 = 2.5
$ws.Range("U186").Value2 = 1.8
$ws.Range("V186").Value2 = 2.05

# Row 187 - now holds match that used to be on row 186 (HamKam vs Molde), with refreshed odds
$ws.Range("B187").Value2 = 7617325
$ws.Range("F187").Value2 = "HamKam"
$ws.Range("G187").Value2 = "Molde"
$ws.Range("K187").Value2 = 4.8
$ws.Range("L187").Value2 = 4.5
$ws.Range("M187").Value2 = 1.571
$ws.Range("N187").Value2 = 5.25
$ws.Range("O187").Value2 = 4.75
$ws.Range("P187").Value2 = 1.5
$ws.Range("Q187").Value2 = 1
$ws.Range("R187").Value2 = 2.02
$ws.Range("S187").Value2 = 1.88
$ws.Range("T187").Value2 = 3
$ws.Range("U187").Value2 = 1.95
$ws.Range("V187").Value2 = 1.9

# Row 188 - odds refresh
$ws.Range("N188").Value2 = 2.45
$ws.Range("P188").Value2 = 2.625
$ws.Range("R188").Value2 = 1.85
$ws.Range("S188").Value2 = 2.05
$ws.Range("U188").Value2 = 1.85
$ws.Range("V188").Value2 = 2

# Row 189 - odds refresh
$ws.Range("N189").Value2 = 1.363
$ws.Range("P189").Value2 = 7.5
$ws.Range("R189").Value2 = 1.95
$ws.Range("S189").Value2 = 1.95
$ws.Range("T189").Value2 = 3.25
$ws.Range("U189").Value2 = 1.9
$ws.Range("V189").Value2 = 1.95
